$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$values = @{
    2  = @{ B = 628;     C = 629.4 }
    3  = @{ B = 3258.8;  C = 3244.05 }
    4  = @{ B = 471.05;  C = 473.95 }
    5  = @{ B = 1569.9;  C = 1599.65 }
    6  = @{ B = 6806.7;  C = 6873.05 }
    7  = @{ B = 191.7;   C = 190.3 }
    8  = @{ B = 263.4;   C = 264.3 }
    9  = @{ B = 48501.35; C = 49182.8 }
    10 = @{ B = 818.3;   C = 820 }
    11 = @{ B = 4745.4;  C = 4770.95 }
    12 = @{ B = 159.5;   C = 160 }
    13 = @{ B = 1353.1;  C = 1358.4 }
    14 = @{ B = 705.3;   C = 681.3 }
    15 = @{ B = 1460.25; C = 1459.35 }
    16 = @{ B = 1038.55; C = 1052.6 }
    17 = @{ B = 646.8;   C = 645.4 }
    18 = @{ B = 2533.65; C = 2552.8 }
    19 = @{ B = 268.45;  C = 269.75 }
    20 = @{ B = 22743.8; C = 22937.25 }
    21 = @{ B = 364.6;   C = 365.4 }
    22 = @{ B = 822.65;  C = 831.15 }
    23 = @{ B = 653.75;  C = 660.7 }
    24 = @{ B = 943.6;   C = 947.55 }
    25 = @{ B = 432.85;  C = 436.9 }
    26 = @{ B = 174.25;  C = 174.9 }
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row].B
    $ws.Range("C$row").Value = $values[$row].C
}
